$p = $ppt.ActivePresentation
$d = $p.Designs.Add()
Write-Output ("count after add=" + $p.Designs.Count)
$d2 = $p.Designs.Item(2)
$d2.Delete()
Write-Output ("count after delete=" + $p.Designs.Count)
